# Weekly refresh of Hortaliza / Agricola del Norte S.A. de Arica - Alcachofa data.
# Applies the updated source rows (dates, variety, volumes, prices) per the
# "Fruta / hortaliza, semanal" data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44412
$ws.Range("H3").Value = 'Symphony'
$ws.Range("J3").Value = 240
# Row 5
$ws.Range("D5").Value = 44377
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 20000
$ws.Range("L5").Value = 21000
$ws.Range("M5").Value = 20333
$ws.Range("P5").Value = 508
# Row 6
$ws.Range("D6").Value = 44377
$ws.Range("H6").Value = 'Symphony'
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 21000
$ws.Range("L6").Value = 22000
$ws.Range("M6").Value = 21500
$ws.Range("P6").Value = 538
# Row 7
$ws.Range("D7").Value = 44435
$ws.Range("H7").Value = 'Madrigal'
$ws.Range("J7").Value = 160
$ws.Range("N7").Value = '$/caja 40 unidades'
$ws.Range("P7").Value = 488
$ws.Range("Q7").Value = 40
# Row 8
$ws.Range("D8").Value = 44356
$ws.Range("K8").Value = 19000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 19500
$ws.Range("P8").Value = 390
# Row 9
$ws.Range("D9").Value = 44363
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 19000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 19500
$ws.Range("P9").Value = 488
# Row 10
$ws.Range("D10").Value = 44426
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 19000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 19500
$ws.Range("P10").Value = 488
# Row 11
$ws.Range("D11").Value = 44384
$ws.Range("H11").Value = 'Madrigal'
$ws.Range("J11").Value = 80
# Row 12
$ws.Range("D12").Value = 44384
$ws.Range("I12").Value = 'Segunda'
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 19000
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 19333
$ws.Range("N12").Value = '$/caja 50 unidades'
$ws.Range("P12").Value = 387
$ws.Range("Q12").Value = 50
# Row 13
$ws.Range("D13").Value = 44384
$ws.Range("H13").Value = 'Symphony'
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 20000
$ws.Range("L13").Value = 21000
$ws.Range("M13").Value = 20400
$ws.Range("P13").Value = 510
# Row 14
$ws.Range("D14").Value = 44483
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14500
$ws.Range("P14").Value = 362
# Row 15
$ws.Range("D15").Value = 44482
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14500
$ws.Range("N15").Value = '$/caja 40 unidades'
$ws.Range("P15").Value = 362
$ws.Range("Q15").Value = 40
# Row 16
$ws.Range("D16").Value = 44370
$ws.Range("H16").Value = 'Argentina(o)'
$ws.Range("J16").Value = 140
$ws.Range("M16").Value = 20429
$ws.Range("N16").Value = '$/caja 50 unidades'
$ws.Range("P16").Value = 409
$ws.Range("Q16").Value = 50
# Row 17
$ws.Range("H17").Value = 'Madrigal'
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 22000
$ws.Range("L17").Value = 23000
$ws.Range("M17").Value = 22500
$ws.Range("N17").Value = '$/caja 40 unidades'
$ws.Range("P17").Value = 562
$ws.Range("Q17").Value = 40
# Row 18
$ws.Range("D18").Value = 44468
$ws.Range("H18").Value = 'Argentina(o)'
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 17000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 17500
$ws.Range("N18").Value = '$/caja 50 unidades'
$ws.Range("P18").Value = 350
$ws.Range("Q18").Value = 50
# Row 19
$ws.Range("D19").Value = 44489
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 13000
$ws.Range("L19").Value = 14000
$ws.Range("M19").Value = 13500
$ws.Range("P19").Value = 338
# Row 20
$ws.Range("D20").Value = 44167
$ws.Range("H20").Value = 'Española'
$ws.Range("J20").Value = 160
$ws.Range("N20").Value = '$/caja 30 unidades'
$ws.Range("O20").Value = 'Región Metropolitana'
$ws.Range("P20").Value = 450
$ws.Range("Q20").Value = 30
# Row 21
$ws.Range("D21").Value = 44160
$ws.Range("J21").Value = 160
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 14500
$ws.Range("P21").Value = 362
# Row 22
$ws.Range("D22").Value = 44391
$ws.Range("H22").Value = 'Madrigal'
$ws.Range("J22").Value = 140
# Row 23
$ws.Range("D23").Value = 44405
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 21000
$ws.Range("L23").Value = 22000
$ws.Range("M23").Value = 21500
$ws.Range("P23").Value = 538
# Row 24
$ws.Range("D24").Value = 44398
$ws.Range("H24").Value = 'Madrigal'
$ws.Range("J24").Value = 170
$ws.Range("K24").Value = 21000
$ws.Range("L24").Value = 22000
$ws.Range("M24").Value = 21500
$ws.Range("N24").Value = '$/caja 40 unidades'
$ws.Range("O24").Value = 'Región de Coquimbo'
$ws.Range("P24").Value = 538
$ws.Range("Q24").Value = 40
